$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2, shifting the bug-tracker table down
$ws.Rows(2).Insert()

# The freshly inserted row picks up formatting from row 1 by default;
# re-stamp it with the formatting of the (now pushed-down) original row 2
# so the new entry keeps the same cell styles as the rest of the table.
$ws.Range("A3:I3").Copy()
$ws.Range("A2:I2").PasteSpecial(-4122)

# Populate the new bug-tracker entry
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = [DateTime]"2016-11-29"
$ws.Range("E2").Value = "core"
$ws.Range("F2").Value = "反射"
$ws.Range("G2").Value = "防骑Q技能对黑曜石雕像"

$ws.Range("G7").Select()
